$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A18 was stored as text ("71277620"); convert it to a real number,
# matching the numeric phone values used elsewhere in column A.
$ws.Range("A18").Value = 71277620

# Append a new redemption row (row 19) for the repeated redemption.
# Column A keeps the phone number as text (leading apostrophe forces
# text entry instead of Excel's automatic numeric conversion); reset
# the style back to Normal afterwards so no new cell style is created.
$ws.Range("A19").Value = "'71277620"
$ws.Range("A19").Style = "Normal"

$ws.Range("B19").Value = 76
$ws.Range("C19").Value = "2025-08-18T17:10:26"
